$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: lowercase column names + two new columns (branch_num, hire_date) ---
$ws.Range("A1").Value = "employee_id"
$ws.Range("B1").Value = "age"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "title"
$ws.Range("E1").Value = "is_admin"
$ws.Range("F1").Value = "branch_num"
$ws.Range("G1").Value = "hire_date"

# --- New column widths for F, G, H ---
# (Target XML widths are 16.7109375 / 21 / 17.42578125; the host's column-width
#  grid is a 1/6-character pixel grid, so F and H land on the closest
#  achievable tick - only G's width of 21 chars is exactly reachable.)
$ws.Columns.Item(6).ColumnWidth = 15.83
$ws.Columns.Item(7).ColumnWidth = 20.17
$ws.Columns.Item(8).ColumnWidth = 16.67

# --- New data: branch_num (F) and hire_date (G) for each employee row ---
# hire_date values given as Excel serial date numbers to avoid locale-dependent
# "format as you type" auto-formatting before the explicit NumberFormat is applied.
$hireDateSerials = @{
    2  = 36682
    3  = 36683
    4  = 40221
    5  = 38701
    6  = 37257
    7  = 41558
    8  = 42179
    9  = 39902
    10 = 39133
    11 = 37410
    12 = 40442
    13 = 40858
    14 = 38135
    15 = 41356
    16 = 38564
}

for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = 20240601
    $ws.Cells.Item($r, 7).Value = $hireDateSerials[$r]
    $ws.Cells.Item($r, 7).NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
}

# --- Update selection to match the author's final cursor position ---
$ws.Range("G17").Select()
